# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 410
    3  = 1394
    4  = 6928
    5  = 425
    6  = 222
    7  = 4473
    8  = 62
    9  = 36
    11 = 900
    13 = 5433
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
